$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.888.89"
$ws.Range("E2").Value = "  -0.41%  "
$ws.Range("D3").Value = "2.262.82"
$ws.Range("E3").Value = "  -0.36%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "304.93"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.27%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "95.39"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.21%  "
$ws.Range("E7").Value = "  -0.75%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("E9").Value = "  +0.41%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.07"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.64%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0790"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.13%  "
$ws.Range("E12").Value = "  -0.52%  "
$ws.Range("E13").Value = "  -0.25%  "
$ws.Range("D14").Value = "2.614.21"
$ws.Range("E14").Value = "  -0.26%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.36"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.04%  "
$ws.Range("D16").Value = "2.254.80"
$ws.Range("E16").Value = "  -0.60%  "
$ws.Range("E17").Value = "  +0.76%  "
$ws.Range("D18").Value = "41.803.89"
$ws.Range("E18").Value = "  -0.34%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.36"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Value = "0.0₃0900"
$ws.Range("E20").Value = "  -1.92%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.95"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.79%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.97"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.30%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "237.03"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.03%  "
$ws.Range("E24").Value = "  -1.74%  "
$ws.Range("E25").Value = "  -0.78%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "23.63"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.67%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "36.72"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.62%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.11"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.04%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.44"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.54%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "160.25"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.16%  "
$ws.Range("E32").Value = "  -2.61%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.999"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.04%  "
$ws.Range("E34").Value = "  +4.95%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0734"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.38%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "16.96"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.31%  "
$ws.Range("E37").Value = "  -0.07%  "
$ws.Range("E38").Value = "  -0.85%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.81"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.97%  "
$ws.Range("E40").Value = "  -2.57%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.99"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.20%  "
$ws.Range("E42").Value = "  +2.17%  "
$ws.Range("E43").Value = "  -2.19%  "
$ws.Range("E44").Value = "  -0.37%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "18.66"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.65%  "
$ws.Range("E46").Value = "  +0.34%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.86"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.90%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "52.91"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.87%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "72.60"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.06%  "
$ws.Range("E50").Value = "  -1.55%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "91.06"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.97%  "
